$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting: copy style of row 74 col A (s=1, bold+border+center) down to A75:A83
$ws.Range("A74").Copy() | Out-Null
$ws.Range("A75:A83").PasteSpecial(-4122) | Out-Null

# Force text (string) number format on numeric-looking text columns so Excel does not
# silently convert these strings into floating-point numbers (losing trailing zeros/precision).
# Row 74: only the cells that are actually changing to a new numeric-looking text value.
$ws.Range("C74").NumberFormat = "@"
$ws.Range("E74:F74").NumberFormat = "@"
$ws.Range("H74").NumberFormat = "@"
$ws.Range("J74:K74").NumberFormat = "@"
# Rows 75:83 are brand new rows, set text format for all the text-typed columns.
$ws.Range("B75:F83").NumberFormat = "@"
$ws.Range("H75:H83").NumberFormat = "@"
$ws.Range("J75:K83").NumberFormat = "@"

# ---- Update existing row 74 (modified values per diff) ----
$ws.Range("C74").Value = "158.97000000"
$ws.Range("E74").Value = "156.29000000"
$ws.Range("F74").Value = "1042091.03622000"
$ws.Range("H74").Value = "160423270.30371500"
$ws.Range("I74").Value = 255965
$ws.Range("J74").Value = "519763.88357000"
$ws.Range("K74").Value = "80051376.52875040"
$ws.Range("M74").Value = 156.2899999999998
$ws.Range("N74").Value = 157.4400000000001
$ws.Range("O74").Value = 162.6114285714285
$ws.Range("P74").Value = 151.9626666666667
$ws.Range("Q74").Value = 140.2810000000001
$ws.Range("R74").Value = 156.29
$ws.Range("S74").Value = 157.1549333544223
$ws.Range("T74").Value = 155.5792030747713
$ws.Range("U74").Value = 154.3681680793534
$ws.Range("V74").Value = 1.211034995417947
$ws.Range("W74").Value = -2.128848630276904
$ws.Range("X74").Value = 3.339883625694851

# ---- Add new rows 75:83 ----
# Row 75
$ws.Range("A75").Value = 73
$ws.Range("B75").Value = "156.28000000"
$ws.Range("C75").Value = "161.85000000"
$ws.Range("D75").Value = "155.23000000"
$ws.Range("E75").Value = "158.45000000"
$ws.Range("F75").Value = "698030.65089000"
$ws.Range("G75").Value = 1586908799999
$ws.Range("H75").Value = "110710086.10400900"
$ws.Range("I75").Value = 178490
$ws.Range("J75").Value = "343955.51426000"
$ws.Range("K75").Value = "54564372.10918630"
$ws.Range("L75").Value = "2020-04-14 08:00:00"
$ws.Range("M75").Value = 158.4499999999998
$ws.Range("N75").Value = 157.3700000000001
$ws.Range("O75").Value = 161.7285714285714
$ws.Range("P75").Value = 153.7220000000001
$ws.Range("Q75").Value = 141.4366666666668
$ws.Range("R75").Value = 158.45
$ws.Range("S75").Value = 158.0183111181408
$ws.Range("T75").Value = 156.0208660296681
$ws.Range("U75").Value = 154.6715460604184
$ws.Range("V75").Value = 1.349319969249649
$ws.Range("W75").Value = -1.43321486348597
$ws.Range("X75").Value = 2.782534832735619

# Row 76
$ws.Range("A76").Value = 74
$ws.Range("B76").Value = "158.46000000"
$ws.Range("C76").Value = "161.29000000"
$ws.Range("D76").Value = "152.00000000"
$ws.Range("E76").Value = "152.73000000"
$ws.Range("F76").Value = "729119.66505000"
$ws.Range("G76").Value = 1586995199999
$ws.Range("H76").Value = "114640840.11551650"
$ws.Range("I76").Value = 186226
$ws.Range("J76").Value = "347034.48291000"
$ws.Range("K76").Value = "54594699.03576640"
$ws.Range("L76").Value = "2020-04-15 08:00:00"
$ws.Range("M76").Value = 152.7299999999998
$ws.Range("N76").Value = 155.5900000000001
$ws.Range("O76").Value = 158.8171428571428
$ws.Range("P76").Value = 155.0560000000001
$ws.Range("Q76").Value = 142.8270000000001
$ws.Range("R76").Value = 152.73
$ws.Range("S76").Value = 154.4927703727136
$ws.Range("T76").Value = 155.5145771154511
$ws.Range("U76").Value = 154.527278688836
$ws.Range("V76").Value = 0.9872984266150979
$ws.Range("W76").Value = -0.9491121793629917
$ws.Range("X76").Value = 1.93641060597809

# Row 77
$ws.Range("A77").Value = 75
$ws.Range("B77").Value = "152.74000000"
$ws.Range("C77").Value = "174.79000000"
$ws.Range("D77").Value = "148.33000000"
$ws.Range("E77").Value = "172.29000000"
$ws.Range("F77").Value = "1647207.59991000"
$ws.Range("G77").Value = 1587081599999
$ws.Range("H77").Value = "273639160.06938750"
$ws.Range("I77").Value = 389854
$ws.Range("J77").Value = "830355.78324000"
$ws.Range("K77").Value = "137897793.15056740"
$ws.Range("L77").Value = "2020-04-16 08:00:00"
$ws.Range("M77").Value = 172.2899999999998
$ws.Range("N77").Value = 162.5100000000001
$ws.Range("O77").Value = 159.2128571428571
$ws.Range("P77").Value = 157.4960000000001
$ws.Range("Q77").Value = 144.7143333333334
$ws.Range("R77").Value = 172.29
$ws.Range("S77").Value = 166.3575901242378
$ws.Range("T77").Value = 158.0954193105433
$ws.Range("U77").Value = 155.84683967885
$ws.Range("V77").Value = 2.248579631693275
$ws.Range("W77").Value = -0.3095737895646652
$ws.Range("X77").Value = 2.55815342125794

# Row 78
$ws.Range("A78").Value = 76
$ws.Range("B78").Value = "172.31000000"
$ws.Range("C78").Value = "174.96000000"
$ws.Range("D78").Value = "168.31000000"
$ws.Range("E78").Value = "170.69000000"
$ws.Range("F78").Value = "723537.75591000"
$ws.Range("G78").Value = 1587167999999
$ws.Range("H78").Value = "123553456.31967860"
$ws.Range("I78").Value = 191289
$ws.Range("J78").Value = "352005.34180000"
$ws.Range("K78").Value = "60121889.86467380"
$ws.Range("L78").Value = "2020-04-17 08:00:00"
$ws.Range("M78").Value = 170.6899999999998
$ws.Range("N78").Value = 171.4900000000001
$ws.Range("O78").Value = 161.0528571428571
$ws.Range("P78").Value = 159.4566666666667
$ws.Range("Q78").Value = 146.467
$ws.Range("R78").Value = 170.69
$ws.Range("S78").Value = 169.2458633747459
$ws.Range("T78").Value = 160.0330521308449
$ws.Range("U78").Value = 156.9492755881899
$ws.Range("V78").Value = 3.083776542655016
$ws.Range("W78").Value = 0.369096300299315
$ws.Range("X78").Value = 2.714680242355701

# Row 79
$ws.Range("A79").Value = 77
$ws.Range("B79").Value = "170.61000000"
$ws.Range("C79").Value = "189.54000000"
$ws.Range("D79").Value = "170.48000000"
$ws.Range("E79").Value = "187.40000000"
$ws.Range("F79").Value = "1127296.36979000"
$ws.Range("G79").Value = 1587254399999
$ws.Range("H79").Value = "203132211.95640910"
$ws.Range("I79").Value = 286172
$ws.Range("J79").Value = "569708.87703000"
$ws.Range("K79").Value = "102696714.57128610"
$ws.Range("L79").Value = "2020-04-18 08:00:00"
$ws.Range("M79").Value = 187.3999999999998
$ws.Range("N79").Value = 179.0450000000001
$ws.Range("O79").Value = 165.2057142857143
$ws.Range("P79").Value = 162.5293333333334
$ws.Range("Q79").Value = 148.1703333333334
$ws.Range("R79").Value = 187.4
$ws.Range("S79").Value = 181.3486211249153
$ws.Range("T79").Value = 164.2433610364284
$ws.Range("U79").Value = 159.2104731794591
$ws.Range("V79").Value = 5.032887856969325
$ws.Range("W79").Value = 1.301854637383962
$ws.Range("X79").Value = 3.731033219585362

# Row 80
$ws.Range("A80").Value = 78
$ws.Range("B80").Value = "187.40000000"
$ws.Range("C80").Value = "188.35000000"
$ws.Range("D80").Value = "175.75000000"
$ws.Range("E80").Value = "180.03000000"
$ws.Range("F80").Value = "995759.57750000"
$ws.Range("G80").Value = 1587340799999
$ws.Range("H80").Value = "181415813.50269480"
$ws.Range("I80").Value = 270668
$ws.Range("J80").Value = "478833.69271000"
$ws.Range("K80").Value = "87227250.30026600"
$ws.Range("L80").Value = "2020-04-19 08:00:00"
$ws.Range("M80").Value = 180.0299999999998
$ws.Range("N80").Value = 183.7150000000001
$ws.Range("O80").Value = 168.2685714285714
$ws.Range("P80").Value = 164.9246666666668
$ws.Range("Q80").Value = 149.722
$ws.Range("R80").Value = 180.03
$ws.Range("S80").Value = 180.4695403749718
$ws.Range("T80").Value = 166.6720792299954
$ws.Range("U80").Value = 160.7561975166961
$ws.Range("V80").Value = 5.915881713299285
$ws.Range("W80").Value = 2.224660072947727
$ws.Range("X80").Value = 3.691221640351558

# Row 81
$ws.Range("A81").Value = 79
$ws.Range("B81").Value = "180.02000000"
$ws.Range("C81").Value = "186.46000000"
$ws.Range("D81").Value = "166.70000000"
$ws.Range("E81").Value = "170.20000000"
$ws.Range("F81").Value = "1547821.52737000"
$ws.Range("G81").Value = 1587427199999
$ws.Range("H81").Value = "273854539.80380100"
$ws.Range("I81").Value = 365589
$ws.Range("J81").Value = "745586.67414000"
$ws.Range("K81").Value = "131918687.66076710"
$ws.Range("L81").Value = "2020-04-20 08:00:00"
$ws.Range("M81").Value = 170.1999999999998
$ws.Range("N81").Value = 175.1150000000001
$ws.Range("O81").Value = 170.2557142857143
$ws.Range("P81").Value = 166.7593333333334
$ws.Range("Q81").Value = 150.9733333333334
$ws.Range("R81").Value = 170.2
$ws.Range("S81").Value = 173.6231801249906
$ws.Range("T81").Value = 167.2148371237567
$ws.Range("U81").Value = 161.4572238087152
$ws.Range("V81").Value = 5.757613315041425
$ws.Range("W81").Value = 2.931250733850843
$ws.Range("X81").Value = 2.826362581190582

# Row 82
$ws.Range("A82").Value = 80
$ws.Range("B82").Value = "170.21000000"
$ws.Range("C82").Value = "174.70000000"
$ws.Range("D82").Value = "168.00000000"
$ws.Range("E82").Value = "170.74000000"
$ws.Range("F82").Value = "906381.65398000"
$ws.Range("G82").Value = 1587513599999
$ws.Range("H82").Value = "155292048.43579020"
$ws.Range("I82").Value = 228262
$ws.Range("J82").Value = "454400.78888000"
$ws.Range("K82").Value = "77856857.48694720"
$ws.Range("L82").Value = "2020-04-21 08:00:00"
$ws.Range("M82").Value = 170.7399999999998
$ws.Range("N82").Value = 170.4700000000001
$ws.Range("O82").Value = 172.0114285714286
$ws.Range("P82").Value = 166.7206666666667
$ws.Range("Q82").Value = 152.5873333333334
$ws.Range("R82").Value = 170.74
$ws.Range("S82").Value = 171.7010600416635
$ws.Range("T82").Value = 167.7571705944951
$ws.Range("U82").Value = 162.1461885365169
$ws.Range("V82").Value = 5.610982057978191
$ws.Range("W82").Value = 3.467197006251793
$ws.Range("X82").Value = 2.143785051726398

# Row 83
$ws.Range("A83").Value = 81
$ws.Range("B83").Value = "170.73000000"
$ws.Range("C83").Value = "184.00000000"
$ws.Range("D83").Value = "169.78000000"
$ws.Range("E83").Value = "181.38000000"
$ws.Range("F83").Value = "812086.87602000"
$ws.Range("G83").Value = 1587599999999
$ws.Range("H83").Value = "143959338.71980660"
$ws.Range("I83").Value = 196518
$ws.Range("J83").Value = "399591.00926000"
$ws.Range("K83").Value = "70844470.39062520"
$ws.Range("L83").Value = "2020-04-22 08:00:00"
$ws.Range("M83").Value = 181.3799999999998
$ws.Range("N83").Value = 176.0600000000001
$ws.Range("O83").Value = 176.1042857142857
$ws.Range("P83").Value = 167.8373333333334
$ws.Range("Q83").Value = 154.1026666666667
$ws.Range("R83").Value = 181.38
$ws.Range("S83").Value = 178.1536866805545
$ws.Range("T83").Value = 169.8529928591669
$ws.Range("U83").Value = 163.5735081370333
$ws.Range("V83").Value = 6.279484722133674
$ws.Range("W83").Value = 4.029654555788338
$ws.Range("X83").Value = 2.249830166345336
